# Update the "想去人数" (F column) counts across the "展览", "演出" and
# "全部类型" worksheets, per the commit's regenerated data.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 (Exhibition) sheet -------------------------------------------------
$wsExhibit.Range("F2").Value  = 914
$wsExhibit.Range("F3").Value  = 511
$wsExhibit.Range("F4").Value  = 511
$wsExhibit.Range("F5").Value  = 792
$wsExhibit.Range("F8").Value  = 813
$wsExhibit.Range("F9").Value  = 462
$wsExhibit.Range("F10").Value = 616
$wsExhibit.Range("F11").Value = 171
$wsExhibit.Range("F13").Value = 30
$wsExhibit.Range("F16").Value = 1540
$wsExhibit.Range("F17").Value = 190
$wsExhibit.Range("F19").Value = 462
$wsExhibit.Range("F20").Value = 57
$wsExhibit.Range("F21").Value = 388
$wsExhibit.Range("F24").Value = 21
$wsExhibit.Range("F25").Value = 210
$wsExhibit.Range("F26").Value = 712
$wsExhibit.Range("F28").Value = 1384
$wsExhibit.Range("F29").Value = 137

# --- 演出 (Show) sheet -------------------------------------------------------
$wsShow.Range("F7").Value = 269

# --- 全部类型 (All types) sheet ----------------------------------------------
$wsAll.Range("F3").Value  = 914
$wsAll.Range("F4").Value  = 511
$wsAll.Range("F5").Value  = 511
$wsAll.Range("F6").Value  = 792
$wsAll.Range("F9").Value  = 813
$wsAll.Range("F12").Value = 462
$wsAll.Range("F13").Value = 616
$wsAll.Range("F15").Value = 171
$wsAll.Range("F17").Value = 30
$wsAll.Range("F20").Value = 1540
$wsAll.Range("F22").Value = 190
$wsAll.Range("F24").Value = 462
$wsAll.Range("F25").Value = 57
$wsAll.Range("F26").Value = 388
$wsAll.Range("F29").Value = 269
$wsAll.Range("F36").Value = 21
$wsAll.Range("F37").Value = 210
$wsAll.Range("F38").Value = 712
$wsAll.Range("F40").Value = 1384
$wsAll.Range("F41").Value = 137
